$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 29400
$ws.Range("J26").Value = 29400
$ws.Range("L26").Value = 29400
$ws.Range("N26").Value = -30088
$ws.Range("H40").Value = 1749.9166
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 1499.9
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 1499.9
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -1849.9
$ws.Range("H111").Value = 2122.5
$ws.Range("I111").Value = 1750
$ws.Range("J111").Value = 2495
$ws.Range("K111").Value = 5250
$ws.Range("L111").Value = 7485
$ws.Range("M111").Value = -2183
$ws.Range("N111").Value = -13619
$ws.Range("H112").Value = 4689.6665
$ws.Range("I112").Value = 35450
$ws.Range("J112").Value = 1451.7368
$ws.Range("K112").Value = 106350
$ws.Range("L112").Value = 4355.2104
$ws.Range("M112").Value = -105242
$ws.Range("N112").Value = -6571.2104
$ws.Range("H118").Value = 859.0769
$ws.Range("I118").Value = 287.5
$ws.Range("K118").Value = 862.5
$ws.Range("M118").Value = 794.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 53649.25
$ws.Range("I28").Value = 27500
$ws.Range("J28").Value = 79798.5
$ws.Range("K28").Value = 27500
$ws.Range("L28").Value = 79798.5
$ws.Range("M28").Value = -27308
$ws.Range("N28").Value = -80182.5
$ws.Range("H45").Value = 1579.5532
$ws.Range("I45").Value = 1515
$ws.Range("K45").Value = 1515
$ws.Range("M45").Value = -1138
$ws.Range("H99").Value = 53649.25
$ws.Range("I99").Value = 27500
$ws.Range("J99").Value = 79798.5
$ws.Range("K99").Value = 27500
$ws.Range("L99").Value = 79798.5
$ws.Range("M99").Value = -24505
$ws.Range("N99").Value = -85788.5
$ws.Range("H122").Value = 125000000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 793.9167
$ws.Range("I107").Value = 296.20514
$ws.Range("J107").Value = 1718.238
$ws.Range("K107").Value = 888.6154199999999
$ws.Range("L107").Value = 5154.714
$ws.Range("M107").Value = 1031.38458
$ws.Range("N107").Value = -8994.714

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 29900
$ws.Range("J117").Value = 29900
$ws.Range("L117").Value = 29900
$ws.Range("N117").Value = -36784
$ws.Range("H122").Value = 23799.8
$ws.Range("I122").Value = 23799.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 71399.39999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -68949.39999999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2885.6
$ws.Range("I126").Value = 1968
$ws.Range("J126").Value = 3636.3635
$ws.Range("K126").Value = 5904
$ws.Range("L126").Value = 10909.0905
$ws.Range("M126").Value = -3434
$ws.Range("N126").Value = -15849.0905
$ws.Range("H133").Value = 60312
$ws.Range("J133").Value = 60312
$ws.Range("L133").Value = 60312
$ws.Range("N133").Value = -70432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 994.5
$ws.Range("I4").Value = 994.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 994.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -881.5
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 15000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -15590
$ws.Range("H28").Value = 994.5
$ws.Range("I28").Value = 994.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 994.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -762.5
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 994.5
$ws.Range("I37").Value = 994.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 994.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -887.5
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 3460
$ws.Range("I40").Value = 3460
$ws.Range("K40").Value = 3460
$ws.Range("M40").Value = -3324
$ws.Range("H46").Value = 1010
$ws.Range("I46").Value = 775
$ws.Range("K46").Value = 775
$ws.Range("M46").Value = -587
$ws.Range("H61").Value = 40175.5
$ws.Range("I61").Value = 40175.5
$ws.Range("K61").Value = 40175.5
$ws.Range("M61").Value = -39973.5
$ws.Range("H93").Value = 4317
$ws.Range("I93").Value = 4707.2354
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 4707.2354
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -3459.2354
$ws.Range("N93").Value = -3496
$ws.Range("H100").Value = 4495
$ws.Range("I100").Value = 3678.9473
$ws.Range("J100").Value = 20000
$ws.Range("K100").Value = 3678.9473
$ws.Range("L100").Value = 20000
$ws.Range("M100").Value = -3137.9473
$ws.Range("N100").Value = -21082
$ws.Range("H113").Value = 40175.5
$ws.Range("I113").Value = 40175.5
$ws.Range("K113").Value = 40175.5
$ws.Range("M113").Value = -38005.5
$ws.Range("H122").Value = 6248.276
$ws.Range("I122").Value = 6341.6665
$ws.Range("J122").Value = 5800
$ws.Range("K122").Value = 19024.9995
$ws.Range("L122").Value = 17400
$ws.Range("M122").Value = -16574.9995
$ws.Range("N122").Value = -22300
$ws.Range("H132").Value = 1873.238
$ws.Range("I132").Value = 1338.2069
$ws.Range("J132").Value = 3066.7693
$ws.Range("K132").Value = 4014.620699999999
$ws.Range("L132").Value = 9200.3079
$ws.Range("M132").Value = -1484.620699999999
$ws.Range("N132").Value = -14260.3079
$ws.Range("H138").Value = 69940
$ws.Range("J138").Value = 69940
$ws.Range("L138").Value = 69940
$ws.Range("N138").Value = -80220

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 49941.8
$ws.Range("J124").Value = 49941.8
$ws.Range("L124").Value = 49941.8
$ws.Range("N124").Value = -59761.8

